# Update TPM-derived values on the active worksheet (Col1a2-Itgb3 LR pair).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.419591
$ws.Range("H2").Value = 4.258773
$ws.Range("I2").Value = 0.001848767113890483
$ws.Range("J2").Value = 0.001848767113890483
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.411950165896
$ws.Range("R2").Value = 3.707551493064
$ws.Range("S2").Value = 0.00006340142672756376
$ws.Range("T2").Value = 0.00006340142672756376

# Row 3
$ws.Range("G3").Value = 1.419591
$ws.Range("H3").Value = 4.258773
$ws.Range("I3").Value = 0.001848767113890483
$ws.Range("J3").Value = 0.001848767113890483
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 10.093789813244
$ws.Range("R3").Value = 90.844108319196
$ws.Range("S3").Value = 0.00155349051469828
$ws.Range("T3").Value = 0.00155349051469828

# Row 4
$ws.Range("G4").Value = 1.419591
$ws.Range("H4").Value = 4.258773
$ws.Range("I4").Value = 0.001848767113890483
$ws.Range("J4").Value = 0.001848767113890483
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 1.506606723133
$ws.Range("R4").Value = 13.559460508197
$ws.Range("S4").Value = 0.0002318751724646395
$ws.Range("T4").Value = 0.0002318751724646395

# Row 5
$ws.Range("I5").Value = 0.9578582377148513
$ws.Range("J5").Value = 0.9578582377148513
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 213.4340539523778
$ws.Range("R5").Value = 1920.9064855714
$ws.Range("S5").Value = 0.03284869057740551
$ws.Range("T5").Value = 0.03284869057740551

# Row 6
$ws.Range("I6").Value = 0.9578582377148513
$ws.Range("J6").Value = 0.9578582377148513
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.8048735157259937
$ws.Range("T6").Value = 0.8048735157259937

# Row 7
$ws.Range("I7").Value = 0.9578582377148513
$ws.Range("J7").Value = 0.9578582377148513
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 780.5827191033695
$ws.Range("R7").Value = 7025.244471930326
$ws.Range("S7").Value = 0.1201360314114521
$ws.Range("T7").Value = 0.1201360314114521

# Row 8
$ws.Range("H8").Value = 92.81792100000001
$ws.Range("I8").Value = 0.04029299517125823
$ws.Range("J8").Value = 0.04029299517125823
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 8.978256872125336
$ws.Range("R8").Value = 80.80431184912801
$ws.Range("S8").Value = 0.001381803777117565
$ws.Range("T8").Value = 0.001381803777117565

# Row 9
$ws.Range("H9").Value = 92.81792100000001
$ws.Range("I9").Value = 0.04029299517125823
$ws.Range("J9").Value = 0.04029299517125823
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.03385758289242332
$ws.Range("T9").Value = 0.03385758289242332

# Row 10
$ws.Range("H10").Value = 92.81792100000001
$ws.Range("I10").Value = 0.04029299517125823
$ws.Range("J10").Value = 0.04029299517125823
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 32.83577307497434
$ws.Range("R10").Value = 295.5219576747691
$ws.Range("S10").Value = 0.005053608501717347
$ws.Range("T10").Value = 0.005053608501717346
